# Weekly update: a new price observation was recorded for
# "Feria Lagunitas de Puerto Montt - Cilantro" and inserted as the new
# row 159 (the data set is kept in reverse-chronological insertion
# order, so every existing row from 159 downward shifts down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 159; this pushes the old
# rows 159..186 down to 160..187 and extends the used range to R187,
# while copying the existing row formatting (e.g. the date style on
# column D) down onto the newly created row, exactly like Excel's
# native "Insert Row" command.
$ws.Rows(159).Insert()

# Populate the new row 159 with the new observation's data.
$newRow = 159
$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item($newRow, 3).Value = 'Los Lagos'
$ws.Cells.Item($newRow, 4).Value = 44504
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112040
$ws.Cells.Item($newRow, 7).Value = 'Cilantro'
$ws.Cells.Item($newRow, 8).Value = 'Sin especificar'
$ws.Cells.Item($newRow, 9).Value = 'Primera'
$ws.Cells.Item($newRow, 10).Value = 150
$ws.Cells.Item($newRow, 11).Value = 9000
$ws.Cells.Item($newRow, 12).Value = 9000
$ws.Cells.Item($newRow, 13).Value = 9000
$ws.Cells.Item($newRow, 14).Value = '$/caja 36 atados'
$ws.Cells.Item($newRow, 15).Value = 'Región Metropolitana'
$ws.Cells.Item($newRow, 16).Value = 250
$ws.Cells.Item($newRow, 17).Value = 36
$ws.Cells.Item($newRow, 18).Value = 'Hortaliza'
